$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the same serial date (45192 -> 2023-09-23) for
# every data row (rows 2 through 482). The commit updates this date to 45202
# (2023-10-03) for all of them, leaving everything else untouched.
$ws.Range("C2:C482").Value = 45202
